$d = $word.ActiveDocument

# 1) Update CSU step text: remove trailing "na página do usuário", add period.
$d.Content.Find.Execute(
    "2- Sistema disponibiliza opção para o usuário de receber notificação na página do usuário",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2- Sistema disponibiliza opção para o usuário de receber notificação.", 2)

# 2) Update step 5 text to mention categoria(s) e subcategoria(s)
$d.Content.Find.Execute(
    "5- Usuário seleciona categorias em que tem interesse e confirma",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5- Usuário seleciona categoria(s) e subcategoria(s) em que tem interesse e confirma", 2)

# 3) Renumber step from 4 to 6 for the "guarda a(s) escolha(s)" line
$d.Content.Find.Execute(
    "4- Sistema guarda a(s) escolha(s) do usuário no banco de dados",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "6- Sistema guarda a(s) escolha(s) do usuário no banco de dados", 2)
